$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    @(2, 1.718974849135873, -0.05066157383224272, 0.6452828968948781),
    @(3, 0.363000991799871, 0.4753304334074437, -0.4791181933660784),
    @(4, -0.8330479243490954, -1.024541131397955, -0.850531605951097),
    @(5, -1.569642914610615, -0.0931295833021941, 1.497903240195856),
    @(6, -0.8567210179775577, 1.614784838779701, 0.6441090578393747),
    @(7, 1.608110189805712, 1.092661623205466, -0.2068012686273616),
    @(8, 0.06094765810494811, -0.5887137739827044, 1.113646459042408),
    @(9, 0.7221143501232178, -0.8408286765292819, -0.5170689656819794),
    @(10, -0.2268550618147382, -0.2738793612729168, 0.3579536068914886),
    @(11, -0.5570260025607773, 0.9372762466130032, 0.7658316593481331)
)

foreach ($row in $values) {
    $r = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
}

$wb.Save()
